$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON): B2, D2, E2 deleted; C2 value updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -0.5911246880189821
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3 (STR): B3:E3 values updated
$ws.Range("B3").Value = -0.75226107008933984
$ws.Range("C3").Value = 0.81958017826096596
$ws.Range("D3").Value = -0.097459496988475572
$ws.Range("E3").Value = 2.3909692343347553

# Update selection to match new reduced range
$ws.Range("B1:E3").Select()
